$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" (Changed) column date serial values from 45174 to 45175
$ws.Range("C2:C5").Value = 45175
